$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.43126688834027505,
    -0.0099999993685067068,
    -0.0089999993537706047,
    -0.011999999821876273,
    -0.044589316932561296,
    -0.0059999993372272797,
    -0.019999999224152631,
    -0.019999999219708187,
    -0.0059999993275541286,
    -0.0059999993251267369,
    -0.0044999993370069546,
    -0.0059999993250370309,
    -0.0059999993251258488,
    -0.011999999277820805,
    0.037061167062650391,
    -0.0059999993247896732,
    -0.0059999993219834735,
    -0.0089999992975462462,
    -0.0089999993748626217,
    -0.07549506592239652,
    -0.0089999993528087074,
    -0.0089999993522216215,
    -0.0089999993479166207,
    -0.041999999075346572,
    -0.041999999070143623,
    -0.0059999993353336833,
    -0.0059999993328014867,
    -0.0059999993213324387,
    -0.011999999266130601,
    -0.0054036160748491824,
    -0.014999999231822159,
    -0.02099999918296902,
    -0.005999999302053638
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
